$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Top Gainers")

$ws.Range("C3").Value = 13.6037
$ws.Range("D3").Value = 9.6327
$ws.Range("E3").Value = 5.1932

$ws.Range("C4").Value = 11.625
$ws.Range("D4").Value = 16.1094
$ws.Range("E4").Value = 16.8684

$ws.Range("C5").Value = 11.5515
$ws.Range("D5").Value = 8.641299999999999
$ws.Range("E5").Value = 9.0909

$ws.Range("C7").Value = 9.4682
$ws.Range("D7").Value = 7.7482
$ws.Range("E7").Value = 12.0464

$ws.Range("C8").Value = 7.7506
$ws.Range("D8").Value = 11.7413
$ws.Range("E8").Value = 14.1315

$ws.Range("B10").Value = "RPOWER"
$ws.Range("C10").Value = 6.7079
$ws.Range("D10").Value = 3.1195
$ws.Range("E10").Value = 5.0721

$ws.Range("B11").Value = "SAIL"
$ws.Range("C11").Value = 6.5678
$ws.Range("D11").Value = 8.7904
$ws.Range("E11").Value = 4.7293

$ws.Range("B12").Value = "WALCHANNAG"
$ws.Range("C12").Value = 6.3459
$ws.Range("D12").Value = 3.8899
$ws.Range("E12").Value = -5.285

$ws.Range("B13").Value = "BUTTERFLY"
$ws.Range("C13").Value = 6.2869
$ws.Range("D13").Value = 9.161
$ws.Range("E13").Value = 11.7278

$ws.Range("B14").Value = "M&MFIN"
$ws.Range("C14").Value = 6.2854
$ws.Range("D14").Value = 6.9093
$ws.Range("E14").Value = 15.7828

$ws.Range("B15").Value = "UTKARSHBNK"
$ws.Range("C15").Value = 6.2628
$ws.Range("D15").Value = -4.6522
$ws.Range("E15").Value = -1.3346

$ws.Range("B16").Value = "ABDL"
$ws.Range("C16").Value = 6.2265
$ws.Range("D16").Value = 5.1185
$ws.Range("E16").Value = 28.0284

$ws.Range("C17").Value = 6.1383
$ws.Range("D17").Value = 3.6495
$ws.Range("E17").Value = 32.2061

$ws.Range("B18").Value = "FISCHER"
$ws.Range("C18").Value = 6.0598
$ws.Range("D18").Value = 11.0589
$ws.Range("E18").Value = 4.2252

$ws.Range("B19").Value = "ADANIENSOL"
$ws.Range("C19").Value = 5.9053
$ws.Range("D19").Value = 3.3201
$ws.Range("E19").Value = 11.8487

$ws.Range("B20").Value = "CELLO"
$ws.Range("C20").Value = 5.491
$ws.Range("D20").Value = 4.3285
$ws.Range("E20").Value = 14.2361

$ws.Range("B21").Value = "IOC"
$ws.Range("C21").Value = 5.4815
$ws.Range("D21").Value = 8.3926
$ws.Range("E21").Value = 8.8123

$ws.Range("C22").Value = 5.4633
$ws.Range("D22").Value = 11.8789
$ws.Range("E22").Value = 12.0803

$ws.Range("B23").Value = "EPACKPEB"
$ws.Range("C23").Value = 5.4253
$ws.Range("D23").Value = -1.091
$ws.Range("E23").Value = "N/A"

$ws.Range("B24").Value = "ABREL"
$ws.Range("C24").Value = 5.2815
$ws.Range("D24").Value = 6.0449
$ws.Range("E24").Value = 5.6055

$ws.Range("C25").Value = 5.2435
$ws.Range("D25").Value = 3.7166
$ws.Range("E25").Value = 6.2867

$ws.Range("C26").Value = 5.1991
$ws.Range("D26").Value = -2.1774
$ws.Range("E26").Value = 17.9031

$ws.Range("C27").Value = 5.153
$ws.Range("D27").Value = 4.8912
$ws.Range("E27").Value = 4.4299

$ws.Range("C28").Value = 5.0934
$ws.Range("D28").Value = 5.6172
$ws.Range("E28").Value = 11.918

$ws.Range("B29").Value = "MEGASOFT"
$ws.Range("C29").Value = 4.9974
$ws.Range("D29").Value = 15.7588
$ws.Range("E29").Value = 33.5271

$ws.Range("B30").Value = "PROZONER"
$ws.Range("C30").Value = 4.9921
$ws.Range("D30").Value = 15.7468
$ws.Range("E30").Value = 36.095

$ws.Range("B31").Value = "STALLION"
$ws.Range("C31").Value = 4.9914
$ws.Range("D31").Value = -5.2229
$ws.Range("E31").Value = 21.4391

$ws.Range("B32").Value = "INDOTHAI"
$ws.Range("C32").Value = 4.9883
$ws.Range("D32").Value = 4.7163
$ws.Range("E32").Value = 43.9974

$ws.Range("B33").Value = "SURYAROSNI"
$ws.Range("C33").Value = 4.7517
$ws.Range("D33").Value = 11.1405
$ws.Range("E33").Value = 2.7943

$ws.Range("B34").Value = "BAJAJINDEF"
$ws.Range("C34").Value = 4.7085
$ws.Range("D34").Value = 3.6272
$ws.Range("E34").Value = 10.6547

$ws.Range("B35").Value = "AXISCADES"
$ws.Range("C35").Value = 4.5839
$ws.Range("D35").Value = 7.0519
$ws.Range("E35").Value = -2.9548

$ws.Range("B36").Value = "GENUSPOWER"
$ws.Range("C36").Value = 4.4932
$ws.Range("D36").Value = 2.8087
$ws.Range("E36").Value = -0.2419

$ws.Range("C37").Value = 4.4726
$ws.Range("D37").Value = 6.6495
$ws.Range("E37").Value = 14.7876

$ws.Range("C38").Value = 4.3595
$ws.Range("D38").Value = 2.8837
$ws.Range("E38").Value = 5.6441

$ws.Range("C39").Value = 4.3561
$ws.Range("D39").Value = 3.1283
$ws.Range("E39").Value = 3.3356

$ws.Range("B40").Value = "FIVESTAR"
$ws.Range("C40").Value = 4.3385
$ws.Range("D40").Value = 4.3482
$ws.Range("E40").Value = 4.426

$ws.Range("B41").Value = "DATAMATICS"
$ws.Range("C41").Value = 4.2514
$ws.Range("D41").Value = 6.6512
$ws.Range("E41").Value = 15.0137

$ws.Range("B42").Value = "LLOYDSENT"
$ws.Range("C42").Value = 4.2483
$ws.Range("D42").Value = 1.5258
$ws.Range("E42").Value = 10.8974

$ws.Range("B43").Value = "STLTECH"
$ws.Range("C43").Value = 4.1752
$ws.Range("D43").Value = 1.1825
$ws.Range("E43").Value = 7.2747

$ws.Range("B44").Value = "STAR"
$ws.Range("C44").Value = 4.1585
$ws.Range("D44").Value = 4.0881
$ws.Range("E44").Value = 3.3207

$ws.Range("B45").Value = "SUNFLAG"
$ws.Range("C45").Value = 4.1296
$ws.Range("D45").Value = 4.466
$ws.Range("E45").Value = 4.7646

$ws.Range("B46").Value = "GMBREW"
$ws.Range("C46").Value = 4.0899
$ws.Range("D46").Value = -0.348
$ws.Range("E46").Value = 79.3565

$ws.Range("C47").Value = 4.0876
$ws.Range("D47").Value = 9.5245
$ws.Range("E47").Value = 19.8509

$ws.Range("B48").Value = "JKIL"
$ws.Range("C48").Value = 4.0247
$ws.Range("D48").Value = 2.8351
$ws.Range("E48").Value = 1.6485

$ws.Range("B49").Value = "RAJRATAN"
$ws.Range("C49").Value = 4.0145
$ws.Range("D49").Value = 1.5501
$ws.Range("E49").Value = 27.7223

$ws.Range("B50").Value = "SUNDROP"
$ws.Range("C50").Value = 3.9379
$ws.Range("D50").Value = 2.422
$ws.Range("E50").Value = 0.5304

$ws.Range("C51").Value = 3.9266
$ws.Range("D51").Value = 1.9634
$ws.Range("E51").Value = 7.7322

$ws.Range("B52").Value = "VINCOFE"
$ws.Range("C52").Value = 3.9183
$ws.Range("D52").Value = 10.7943
$ws.Range("E52").Value = 9.166600000000001

$ws.Range("B53").Value = "GPPL"
$ws.Range("C53").Value = 3.8933
$ws.Range("D53").Value = 2.8902
$ws.Range("E53").Value = 4.5245

$ws.Range("B54").Value = "SHK"
$ws.Range("C54").Value = 3.8361
$ws.Range("D54").Value = 2.5869
$ws.Range("E54").Value = -1.7414

$ws.Range("B55").Value = "TDPOWERSYS"
$ws.Range("C55").Value = 3.8322
$ws.Range("D55").Value = 4.8754
$ws.Range("E55").Value = 13.9545

$ws.Range("B56").Value = "BLACKBUCK"
$ws.Range("C56").Value = 3.8127
$ws.Range("D56").Value = 2.427
$ws.Range("E56").Value = 8.404199999999999

$ws.Range("B57").Value = "PROSTARM"
$ws.Range("C57").Value = 3.8018
$ws.Range("D57").Value = 0.9227
$ws.Range("E57").Value = -8.0006

$ws.Range("B58").Value = "RECLTD"
$ws.Range("C58").Value = 3.7946
$ws.Range("D58").Value = 2.7778
$ws.Range("E58").Value = 2.7089

$ws.Range("B59").Value = "HCC"
$ws.Range("C59").Value = 3.7803
$ws.Range("D59").Value = 2.6817
$ws.Range("E59").Value = 7.4197

$ws.Range("B60").Value = "DCMSHRIRAM"
$ws.Range("C60").Value = 3.7779
$ws.Range("D60").Value = 10.3536
$ws.Range("E60").Value = 17.7367

$ws.Range("B61").Value = "NBCC"
$ws.Range("C61").Value = 3.7621
$ws.Range("D61").Value = 2.4801
$ws.Range("E61").Value = 6.892

$ws.Range("B62").Value = "HITECHGEAR"
$ws.Range("C62").Value = 3.7113
$ws.Range("D62").Value = 1.0051
$ws.Range("E62").Value = 9.769299999999999

$ws.Range("B63").Value = "TCI"
$ws.Range("C63").Value = 3.7103
$ws.Range("D63").Value = 3.6142
$ws.Range("E63").Value = 4.1141

$ws.Range("B64").Value = "SALASAR"
$ws.Range("C64").Value = 3.5827
$ws.Range("D64").Value = 4.5745
$ws.Range("E64").Value = 10.823

$ws.Range("B65").Value = "SGMART"
$ws.Range("C65").Value = 3.5402
$ws.Range("D65").Value = 7.4974
$ws.Range("E65").Value = 1.8169

$ws.Range("B66").Value = "PRAKASH"
$ws.Range("C66").Value = 3.5225
$ws.Range("D66").Value = 4.4259
$ws.Range("E66").Value = 1.1739

$ws.Range("B67").Value = "INDORAMA"
$ws.Range("C67").Value = 3.4913
$ws.Range("D67").Value = 2.6571
$ws.Range("E67").Value = 13.7823

$ws.Range("B68").Value = "RESPONIND"
$ws.Range("C68").Value = 3.4641
$ws.Range("D68").Value = 7.4713
$ws.Range("E68").Value = 6.2976

$ws.Range("B69").Value = "SWANCORP"
$ws.Range("C69").Value = 3.4614
$ws.Range("D69").Value = 12.9942
$ws.Range("E69").Value = 7.4775

$ws.Range("B70").Value = "RHIM"
$ws.Range("C70").Value = 3.4309
$ws.Range("D70").Value = 3.005
$ws.Range("E70").Value = 4.9558

$ws.Range("B71").Value = "TATVA"
$ws.Range("C71").Value = 3.4023
$ws.Range("D71").Value = 3.4943
$ws.Range("E71").Value = 37.8595

$ws.Range("B72").Value = "ORIENTTECH"
$ws.Range("C72").Value = 3.3956
$ws.Range("D72").Value = 0.1071
$ws.Range("E72").Value = 32.1272

$ws.Range("B73").Value = "PENIND"
$ws.Range("C73").Value = 3.3908
$ws.Range("D73").Value = 2.6711
$ws.Range("E73").Value = 12.0247

$ws.Range("B74").Value = "TMB"
$ws.Range("C74").Value = 3.378
$ws.Range("D74").Value = 7.0754
$ws.Range("E74").Value = 14.2315

$ws.Range("B75").Value = "MAITHANALL"
$ws.Range("C75").Value = 3.3478
$ws.Range("D75").Value = 2.3914
$ws.Range("E75").Value = 1.6649

$ws.Range("B76").Value = "SALZERELEC"
$ws.Range("C76").Value = 3.3082
$ws.Range("D76").Value = 6.1192
$ws.Range("E76").Value = 16.4952
